$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '30.824.40'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  -0.96%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.940.40'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  -0.70%  '

$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  +0.10%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '242.77'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  -1.29%  '

$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  +0.06%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4890'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  -0.10%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2935'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  -0.85%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06916'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  +1.25%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '19.27'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  +0.23%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '105.42'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  -1.16%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.932.30'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  -0.93%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.07720'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  -0.14%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.361'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  -0.92%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.6986'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  -2.09%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '272.41'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  -4.59%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '30.831.13'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  -0.54%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.000007717'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  -0.47%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '13.08'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  -1.00%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '2.196.44'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  +0.44%  '

$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  +0.08%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.524'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  +0.12%  '

$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  +0.12%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '6.546'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  -0.89%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '9.723'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  -2.00%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '166.36'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  -1.47%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '19.59'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  -1.95%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.169'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  -1.79%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.1038'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  -1.50%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.388'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  -3.47%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.570'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  -3.70%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.555'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  -2.37%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.364'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  -2.54%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.04861'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  -2.90%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.7555'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  -0.92%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.155'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  -0.98%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.9999'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  +0.11%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.725'

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.01995'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  -2.39%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.658'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  -1.71%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '6.498'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  +1.08%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '77.21'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  +6.30%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.088'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  -2.84%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.9042'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  +2.68%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.4407'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  -1.87%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '107.71'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  -1.69%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.9987'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  -0.07%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '7.748'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  +3.34%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '991.28'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  +1.00%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.1246'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  -2.95%  '

$ws.Range('B51').NumberFormat = '@'
$ws.Range('B51').Value = 'Elrond'
$ws.Range('C51').NumberFormat = '@'
$ws.Range('C51').Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '36.08'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  +0.06%  '
